$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.145036666666667
$ws.Range("H2").Value = 3.43511
$ws.Range("I2").Value = 0.4953865629219574
$ws.Range("J2").Value = 0.4953865629219574
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 58.23819616223778
$ws.Range("R2").Value = 524.14376546014
$ws.Range("S2").Value = 0.1701313290138397
$ws.Range("T2").Value = 0.1701313290138397
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.145036666666667
$ws.Range("H3").Value = 3.43511
$ws.Range("I3").Value = 0.4953865629219574
$ws.Range("J3").Value = 0.4953865629219574
$ws.Range("M3").Value = 43.683024
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("Q3").Value = 50.01866419088
$ws.Range("R3").Value = 450.16797771792
$ws.Range("S3").Value = 0.1461195980484225
$ws.Range("T3").Value = 0.1461195980484226
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.145036666666667
$ws.Range("H4").Value = 3.43511
$ws.Range("I4").Value = 0.4953865629219574
$ws.Range("J4").Value = 0.4953865629219574
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 41.95827208720778
$ws.Range("R4").Value = 377.62444878487
$ws.Range("S4").Value = 0.1225727626150202
$ws.Range("T4").Value = 0.1225727626150202
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.145036666666667
$ws.Range("H5").Value = 3.43511
$ws.Range("I5").Value = 0.4953865629219574
$ws.Range("J5").Value = 0.4953865629219574
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 19.36221698036111
$ws.Range("R5").Value = 174.25995282325
$ws.Range("S5").Value = 0.05656287324467495
$ws.Range("T5").Value = 0.05656287324467495
$ws.Range("G6").Value = 0.4713496666666666
$ws.Range("I6").Value = 0.2039238551060172
$ws.Range("J6").Value = 0.2039238551060172
$ws.Range("M6").Value = 50.86142466666666
$ws.Range("N6").Value = 152.584274
$ws.Range("O6").Value = 0.3434314568613803
$ws.Range("P6").Value = 0.3434314568613804
$ws.Range("Q6").Value = 23.97351556282511
$ws.Range("R6").Value = 215.761640065426
$ws.Range("S6").Value = 0.07003386664784852
$ws.Range("T6").Value = 0.07003386664784854
$ws.Range("G7").Value = 0.4713496666666666
$ws.Range("I7").Value = 0.2039238551060172
$ws.Range("J7").Value = 0.2039238551060172
$ws.Range("M7").Value = 43.683024
$ws.Range("O7").Value = 0.294960761928139
$ws.Range("P7").Value = 0.294960761928139
$ws.Range("Q7").Value = 20.589978801392
$ws.Range("R7").Value = 185.309809212528
$ws.Range("S7").Value = 0.06014953567739426
$ws.Range("T7").Value = 0.06014953567739427
$ws.Range("G8").Value = 0.4713496666666666
$ws.Range("I8").Value = 0.2039238551060172
$ws.Range("J8").Value = 0.2039238551060172
$ws.Range("M8").Value = 36.64360566666667
$ws.Range("N8").Value = 109.930817
$ws.Range("O8").Value = 0.2474285170192034
$ws.Range("P8").Value = 0.2474285170192035
$ws.Range("Q8").Value = 17.27195131644811
$ws.Range("R8").Value = 155.447561848033
$ws.Range("S8").Value = 0.05045657705372076
$ws.Range("T8").Value = 0.05045657705372077
$ws.Range("G9").Value = 0.4713496666666666
$ws.Range("I9").Value = 0.2039238551060172
$ws.Range("J9").Value = 0.2039238551060172
$ws.Range("M9").Value = 16.90969166666667
$ws.Range("N9").Value = 50.729075
$ws.Range("O9").Value = 0.1141792641912772
$ws.Range("P9").Value = 0.1141792641912772
$ws.Range("Q9").Value = 7.970377530519444
$ws.Range("R9").Value = 71.733397774675
$ws.Range("S9").Value = 0.02328387572705368
$ws.Range("T9").Value = 0.02328387572705368
$ws.Range("G10").Value = 0.6323219999999999
$ws.Range("H10").Value = 1.896966
$ws.Range("I10").Value = 0.2735666300991275
$ws.Range("J10").Value = 0.2735666300991275
$ws.Range("M10").Value = 50.86142466666666
$ws.Range("N10").Value = 152.584274
$ws.Range("O10").Value = 0.3434314568613803
$ws.Range("P10").Value = 0.3434314568613804
$ws.Range("Q10").Value = 32.16079776807599
$ws.Range("R10").Value = 289.447179912684
$ws.Range("S10").Value = 0.09395138632360168
$ws.Range("T10").Value = 0.09395138632360171
$ws.Range("G11").Value = 0.6323219999999999
$ws.Range("H11").Value = 1.896966
$ws.Range("I11").Value = 0.2735666300991275
$ws.Range("J11").Value = 0.2735666300991275
$ws.Range("M11").Value = 43.683024
$ws.Range("O11").Value = 0.294960761928139
$ws.Range("P11").Value = 0.294960761928139
$ws.Range("Q11").Value = 27.62173710172799
$ws.Range("R11").Value = 248.595633915552
$ws.Range("S11").Value = 0.080691421652152
$ws.Range("T11").Value = 0.08069142165215203
$ws.Range("G12").Value = 0.6323219999999999
$ws.Range("H12").Value = 1.896966
$ws.Range("I12").Value = 0.2735666300991275
$ws.Range("J12").Value = 0.2735666300991275
$ws.Range("M12").Value = 36.64360566666667
$ws.Range("N12").Value = 109.930817
$ws.Range("O12").Value = 0.2474285170192034
$ws.Range("P12").Value = 0.2474285170192035
$ws.Range("Q12").Value = 23.170558022358
$ws.Range("R12").Value = 208.535022201222
$ws.Range("S12").Value = 0.06768818559136809
$ws.Range("T12").Value = 0.06768818559136811
$ws.Range("G13").Value = 0.6323219999999999
$ws.Range("H13").Value = 1.896966
$ws.Range("I13").Value = 0.2735666300991275
$ws.Range("J13").Value = 0.2735666300991275
$ws.Range("M13").Value = 16.90969166666667
$ws.Range("N13").Value = 50.729075
$ws.Range("O13").Value = 0.1141792641912772
$ws.Range("P13").Value = 0.1141792641912772
$ws.Range("Q13").Value = 10.69237005405
$ws.Range("R13").Value = 96.23133048645001
$ws.Range("S13").Value = 0.03123563653200568
$ws.Range("T13").Value = 0.03123563653200569
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.062692
$ws.Range("H14").Value = 0.188076
$ws.Range("I14").Value = 0.02712295187289783
$ws.Range("J14").Value = 0.02712295187289783
$ws.Range("M14").Value = 50.86142466666666
$ws.Range("N14").Value = 152.584274
$ws.Range("O14").Value = 0.3434314568613803
$ws.Range("P14").Value = 0.3434314568613804
$ws.Range("Q14").Value = 3.188604435202667
$ws.Range("R14").Value = 28.697439916824
$ws.Range("S14").Value = 0.009314874876090406
$ws.Range("T14").Value = 0.009314874876090408
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.062692
$ws.Range("H15").Value = 0.188076
$ws.Range("I15").Value = 0.02712295187289783
$ws.Range("J15").Value = 0.02712295187289783
$ws.Range("M15").Value = 43.683024
$ws.Range("O15").Value = 0.294960761928139
$ws.Range("P15").Value = 0.294960761928139
$ws.Range("Q15").Value = 2.738576140608
$ws.Range("R15").Value = 24.647185265472
$ws.Range("S15").Value = 0.008000206550170188
$ws.Range("T15").Value = 0.00800020655017019
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.062692
$ws.Range("H16").Value = 0.188076
$ws.Range("I16").Value = 0.02712295187289783
$ws.Range("J16").Value = 0.02712295187289783
$ws.Range("M16").Value = 36.64360566666667
$ws.Range("N16").Value = 109.930817
$ws.Range("O16").Value = 0.2474285170192034
$ws.Range("P16").Value = 0.2474285170192035
$ws.Range("Q16").Value = 2.297260926454666
$ws.Range("R16").Value = 20.675348338092
$ws.Range("S16").Value = 0.006710991759094336
$ws.Range("T16").Value = 0.006710991759094337
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.062692
$ws.Range("H17").Value = 0.188076
$ws.Range("I17").Value = 0.02712295187289783
$ws.Range("J17").Value = 0.02712295187289783
$ws.Range("M17").Value = 16.90969166666667
$ws.Range("N17").Value = 50.729075
$ws.Range("O17").Value = 0.1141792641912772
$ws.Range("P17").Value = 0.1141792641912772
$ws.Range("Q17").Value = 1.060102389966667
$ws.Range("R17").Value = 9.5409215097
$ws.Range("S17").Value = 0.003096878687542898
$ws.Range("T17").Value = 0.003096878687542898
